$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (their matches are being moved up into rows 2-3 with refreshed odds)
$ws.Rows("4:5").Delete()

# Row 2: Argentinian Primera Division / Instituto vs Velez Sarsfield, 22:15:00
$ws.Range("C2").Value = "22:15:00"
$ws.Range("D2").Value = "Instituto"
$ws.Range("E2").Value = "Velez Sarsfield"
$ws.Range("F2").Value = 2.82
$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.84
$ws.Range("K2").Value = 2.88
$ws.Range("L2").Value = 1.89
$ws.Range("M2").Value = 1.2
$ws.Range("N2").Value = 2.04
$ws.Range("O2").Value = 1.93
$ws.Range("P2").Value = 1.31
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 1.1
$ws.Range("S2").Value = 9.800000000000001
$ws.Range("T2").Value = 2.68
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 1.41
$ws.Range("W2").Value = 1.52
$ws.Range("X2").Value = 5.6
$ws.Range("Y2").Value = 7.8
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 90
$ws.Range("AB2").Value = 6.6
$ws.Range("AC2").Value = 7.2
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 80
$ws.Range("AF2").Value = 16.5
$ws.Range("AG2").Value = 17.5
$ws.Range("AH2").Value = 40
$ws.Range("AI2").Value = 160
$ws.Range("AJ2").Value = 60
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 150
$ws.Range("AM2").Value = 540
$ws.Range("AN2").Value = 100
$ws.Range("AO2").Value = 160

# Row 3: Argentinian Primera Division / Central Cordoba (SdE) vs Gimnasia Mendoza, 22:15:00
$ws.Range("C3").Value = "22:15:00"
$ws.Range("D3").Value = "Central Cordoba (SdE)"
$ws.Range("E3").Value = "Gimnasia Mendoza"
$ws.Range("F3").Value = 6.2
$ws.Range("G3").Value = 7.4
$ws.Range("H3").Value = 1.66
$ws.Range("I3").Value = 1.73
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 3.95
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 6.2
$ws.Range("O3").Value = 1.15
$ws.Range("P3").Value = 1.35
$ws.Range("Q3").Value = 3.45
$ws.Range("R3").Value = 1.13
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 1.05
$ws.Range("U3").Value = 2.58
$ws.Range("V3").Value = 2.02
$ws.Range("W3").Value = 1.12
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 7.6
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

$wb.Save()
